# Adds loca, cleans up button prefabs
# Appends new localization rows for "Back" / "PlayAgain" to the Localization sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# Row 16: Back / Back / Zurück / (Polish not yet translated)
$ws.Cells.Item(16, 1).Value = "Back"
$ws.Cells.Item(16, 2).Value = "Back"
$ws.Cells.Item(16, 3).Value = "Zurück"
$ws.Cells.Item(16, 4).Formula = "=""""" 

# Row 17: PlayAgain / Play Again / Nochmal Spielen / (Polish not yet translated)
$ws.Cells.Item(17, 1).Value = "PlayAgain"
$ws.Cells.Item(17, 2).Value = "Play Again"
$ws.Cells.Item(17, 3).Value = "Nochmal Spielen"
$ws.Cells.Item(17, 4).Formula = "=""""" 

$ws.Activate()
$ws.Rows.Item(17).Select()
